# Daily update at 8 AM UTC
# Appends the next day's row (2025-05-08) to the "Wins Over Time" tracker
# and moves the "last row" date formatting down from row 44 to the new row 45.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 44 was previously the last row and used the short date format
# (YYYY-MM-DD). Now that it's no longer the last row, give it the same
# date/time format used by the rest of the data rows above it.
$ws.Range("A44").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Append the new day's data as row 45.
$ws.Range("A45").Value = 45785
$ws.Range("B45").Value = 182
$ws.Range("C45").Value = 195
$ws.Range("D45").Value = 185

# The newly-appended last row takes on the short date format previously
# used by row 44.
$ws.Range("A45").NumberFormat = "YYYY-MM-DD"
